$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-All "228.84" "223.29"
Replace-All "11.10" "5.55"
Replace-All "87.84" "82.29"
Replace-All "232.24" "226.70"
Replace-All "89.22" "83.67"
Replace-All "234.49" "228.94"
Replace-All "96.11" "90.56"
Replace-All "235.00" "229.45"
Replace-All "89.64" "84.09"
Replace-All "235.85" "230.30"
Replace-All "90.12" "84.57"
Replace-All "236.88" "231.33"
Replace-All "96.50" "90.95"
Replace-All "236.96" "231.42"
Replace-All "96.66" "91.11"
Replace-All "237.38" "231.83"
Replace-All "102.74" "97.19"
Replace-All "238.54" "232.99"
Replace-All "91.08" "85.53"
Replace-All "239.02" "233.48"
Replace-All "103.60" "98.05"
Replace-All "61.6%" "63.1%"
Replace-All "38.4%" "36.9%"
Replace-All "510.06" "504.52"
Replace-All "286.98" "281.43"
Replace-All "274.38" "268.83"
